# The authored edit swaps the raw contents of the two theme parts that ship
# with this deck: ppt/theme/theme1.xml ("Office Theme" colours, only used by
# the Notes Master) and ppt/theme/theme2.xml ("Integral" colours, the theme
# actually applied to the slide master / slides). Both themes already share
# an identical font scheme (Arial-based "Office") and an identical fmtScheme
# (fill/line/effect styles are all expressed relative to `phClr`), so the
# only real content difference between the two parts is their 12-colour
# `<a:clrScheme>` palette. Swapping the two files is therefore equivalent,
# from the rendered deck's point of view, to swapping the 12 theme colours
# that are actually applied to the presentation.
#
# PowerPoint's object model only exposes a single "live" theme colour
# scheme (reachable from the slide/master/handout-master/notes-master -
# they all resolve to the one theme that backs the slide master), so we
# push the colours that originally lived in theme1.xml ("Office") onto
# that live ThemeColorScheme. That is the same effect the authored edit
# has on everything PowerPoint actually renders from the slide side.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (scheme slot, original theme1.xml "Office" RGB hex)
#  1 dk1      000000
#  2 lt1      FFFFFF
#  3 dk2      44546A
#  4 lt2      E7E6E6
#  5 accent1  5B9BD5
#  6 accent2  ED7D31
#  7 accent3  A5A5A5
#  8 accent4  FFC000
#  9 accent5  4472C4
# 10 accent6  70AD47
# 11 hlink    0563C1
# 12 folHlink 954F72
$officeHex = @(
  "000000",
  "FFFFFF",
  "44546A",
  "E7E6E6",
  "5B9BD5",
  "ED7D31",
  "A5A5A5",
  "FFC000",
  "4472C4",
  "70AD47",
  "0563C1",
  "954F72"
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's RGB colour integers are packed 0x00BBGGRR.
    $tcs.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
